$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-word the "limitations" paragraphs (B12:B16) and give them the new
#    bold / accent-coloured look. Clear first so no stale formatting from
#    the previous (unrelated) text survives.
# ---------------------------------------------------------------------------
$ws.Range("B12:B17").Clear()

$ws.Range("B12").Value = "You may consider the following before building your own Microsoft Excel template:"
$ws.Range("B13").Value = "Merging Microsoft Excel templates with OpenTBS has several limitations because of the OpenXML format for Excel."
$ws.Range("B14").Value = "* Formulas won't work because OpenTBS needs to convert cell positions from aboslute to relative in order to have a constistent merged sheet."
$ws.Range("B15").Value = "* Formulas may also make troubles because they are saved twice in the sheet:  one for the expression, and one for the instant result."
$ws.Range("B16").Value = "* Changing picture (using ope=changepic)  because drawing information are saved in another XML sub-file."

$limits = $ws.Range("B12:B16")
$limits.Font.Bold = $true
$limits.Font.ThemeColor = 6
$limits.Font.TintAndShade = -0.249977111117893

# ---------------------------------------------------------------------------
# 2. New "Example #1" sub-title (old B18/B19 text is removed).
# ---------------------------------------------------------------------------
$ws.Range("B18:B19").Clear()
$ws.Range("B18").Value = "Example #1: merging data with rows"
$ws.Range("B18").Font.Bold = $true
$ws.Range("B18").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 3. New little table: header row (B20:D20) + sample data row (B21:D21).
# ---------------------------------------------------------------------------
$ws.Range("B20:D21").Clear()

$ws.Range("B20").Value = "First Name"
$ws.Range("C20").Value = "Name"
$ws.Range("D20").Value = "Membership number"

$ws.Range("B21").Value = "[a.firstname;block=row]"
$ws.Range("C21").Value = "[a.name]"
$ws.Range("D21").Value = "[a.number]"

$header = $ws.Range("B20:D20")
$header.Font.Bold = $false
$header.Interior.ThemeColor = 2
$header.Interior.TintAndShade = -0.14999847407452621
$header.Borders.LineStyle = 1

$dataRow = $ws.Range("B21:D21")
$dataRow.Font.Bold = $false
$dataRow.Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4. Column widths for the new table.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 16.2
$ws.Columns.Item(3).ColumnWidth = 12.2
$ws.Columns.Item(4).ColumnWidth = 19.2

# ---------------------------------------------------------------------------
# 5. Selection / active cell moves to B16 after the edits.
# ---------------------------------------------------------------------------
$ws.Range("B16").Select() | Out-Null
